$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F4").Value = 289
$ws1.Range("F5").Value = 285
$ws1.Range("F6").Value = 385
$ws1.Range("F7").Value = 840
$ws1.Range("F8").Value = 31
$ws1.Range("F9").Value = 492
$ws1.Range("F12").Value = 13
$ws1.Range("F17").Value = 6543
$ws1.Range("F21").Value = 7478
$ws1.Range("F22").Value = 36
$ws1.Range("F24").Value = 3364
$ws1.Range("G24").Value = 63
$ws1.Range("F26").Value = 1143
$ws1.Range("F27").Value = 873
$ws1.Range("F28").Value = 4509
$ws1.Range("F29").Value = 10
$ws1.Range("F31").Value = 62
$ws1.Range("F32").Value = 195
$ws1.Range("F33").Value = 189
$ws1.Range("I33").Value = "//i1.hdslb.com/bfs/openplatform/202402/P1YCG3MT1708329896103.jpeg"
$ws1.Range("F34").Value = 1565
$ws1.Range("F35").Value = 5
$ws1.Range("F38").Value = 16
$ws1.Range("F39").Value = 1155
$ws1.Range("F40").Value = 1658
$ws1.Range("F41").Value = 2124
$ws3.Range("F3").Value = 1210
$ws4.Range("F4").Value = 1210
$ws4.Range("F7").Value = 285
$ws4.Range("F8").Value = 385
$ws4.Range("F9").Value = 840
$ws4.Range("F10").Value = 31
$ws4.Range("F11").Value = 492
$ws4.Range("F15").Value = 14
$ws4.Range("F21").Value = 6543
$ws4.Range("F25").Value = 7478
$ws4.Range("F26").Value = 36
$ws4.Range("F28").Value = 3364
$ws4.Range("G28").Value = 63
$ws4.Range("F30").Value = 1143
$ws4.Range("F31").Value = 873
$ws4.Range("F32").Value = 4509
$ws4.Range("F33").Value = 10
$ws4.Range("F35").Value = 62
$ws4.Range("F37").Value = 195
$ws4.Range("F38").Value = 189
$ws4.Range("I38").Value = "//i1.hdslb.com/bfs/openplatform/202402/P1YCG3MT1708329896103.jpeg"
$ws4.Range("F39").Value = 1565
$ws4.Range("F40").Value = 5
$ws4.Range("F43").Value = 16
$ws4.Range("F44").Value = 1155
$ws4.Range("F45").Value = 1658
$ws4.Range("F47").Value = 2124
